$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.560697
$ws.Range("N2").Value = 1.682091
$ws.Range("O2").Value = 0.05296871374736462
$ws.Range("P2").Value = 0.05296871374736462
$ws.Range("Q2").Value = 0.034043839749
$ws.Range("R2").Value = 0.306394557741
$ws.Range("S2").Value = 0.05296871374736462
$ws.Range("T2").Value = 0.05296871374736462

# Row 3 (FAPs -> FAPs)
$ws.Range("O3").Value = 0.3900427732288309
$ws.Range("P3").Value = 0.3900427732288309
$ws.Range("S3").Value = 0.3900427732288309
$ws.Range("T3").Value = 0.3900427732288309

# Row 4 (FAPs -> MuSCs)
$ws.Range("O4").Value = 0.5569885130238045
$ws.Range("P4").Value = 0.5569885130238045
$ws.Range("S4").Value = 0.5569885130238045
$ws.Range("T4").Value = 0.5569885130238045
